$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: insert two new rows at the top of the data (row 4),
# pushing all existing records down by two rows, then populate the new
# rows with this week's figures.
$ws.Rows("4:5").Insert()

# Row 4 - Primera
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Femacal de La Calera"
$ws.Range("C4").Value = "Coquimbo"
$ws.Range("D4").Value = 44991
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 100112043
$ws.Range("G4").Value = "Pepino dulce"
$ws.Range("H4").Value = "Cultivar IV Región"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 75
$ws.Range("K4").Value = 13000
$ws.Range("L4").Value = 13000
$ws.Range("M4").Value = 13000
$ws.Range("N4").Value = "$/bandeja 18 kilos"
$ws.Range("O4").Value = "Provincia de Limarí"
$ws.Range("P4").Value = 722
$ws.Range("Q4").Value = 18
$ws.Range("R4").Value = "Hortaliza"

# Row 5 - Segunda
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Femacal de La Calera"
$ws.Range("C5").Value = "Coquimbo"
$ws.Range("D5").Value = 44991
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 100112043
$ws.Range("G5").Value = "Pepino dulce"
$ws.Range("H5").Value = "Cultivar IV Región"
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 56
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 9000
$ws.Range("N5").Value = "$/bandeja 18 kilos"
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 500
$ws.Range("Q5").Value = 18
$ws.Range("R5").Value = "Hortaliza"
